# Update Import Feature And Add Order Column
#
# Adds a new "Order" column (H) to the data-import template header row,
# extends the title merge/style to cover it, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Order" header cell (H2), matching the text-format style already
#     used by the other header cells (A2:G2). ---
$ws.Range("H2").Value = "Order"
$ws.Range("H2").NumberFormat = "@"

# --- Give H1 (part of the title band) the same centered text-format style
#     already used by A1:G1. ---
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").HorizontalAlignment = -4108

# --- Extend the title merge from A1:G1 to A1:H1. ---
$ws.Range("A1:H1").Merge()

# --- Size the new column similarly to the other narrow columns. ---
$ws.Columns(8).ColumnWidth = 5

# --- Move the active selection to F4. ---
$null = $ws.Range("F4").Select()
